# Re-applies the updated coinranking.com snapshot to the Sheet1 table (rows 2-51).
# Values are written per-column to match the source data types exactly:
#  - Coin / Link (B, C) and Volume(1h) (E) are plain text.
#  - Price (D) is text too (the sheet never stores it as a number), so numeric-
#    looking prices are written with a leading apostrophe to stop Excel from
#    auto-converting them to the Number type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.477.32'
$ws.Range('E2').Value = '  +2.68%  '

$ws.Range('D3').Value = '2.354.10'
$ws.Range('E3').Value = '  +6.10%  '

$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').Value = "'312.38"
$ws.Range('E5').Value = '  +5.47%  '

$ws.Range('D6').Value = "'109.82"
$ws.Range('E6').Value = '  +1.95%  '

$ws.Range('D7').Value = "'0.644"
$ws.Range('E7').Value = '  +4.04%  '

$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('D9').Value = "'0.636"
$ws.Range('E9').Value = '  +6.36%  '

$ws.Range('D10').Value = "'43.32"
$ws.Range('E10').Value = '  -0.32%  '

$ws.Range('D11').Value = "'0.0939"
$ws.Range('E11').Value = '  +2.93%  '

$ws.Range('D12').Value = "'8.88"
$ws.Range('E12').Value = '  +1.58%  '

$ws.Range('E13').Value = '  +4.71%  '

$ws.Range('E14').Value = '  +2.42%  '

$ws.Range('D15').Value = "'16.39"
$ws.Range('E15').Value = '  +8.77%  '

$ws.Range('D16').Value = '2.707.37'
$ws.Range('E16').Value = '  +6.11%  '

$ws.Range('D17').Value = '2.415.77'
$ws.Range('E17').Value = '  +8.04%  '

$ws.Range('D18').Value = '43.447.52'
$ws.Range('E18').Value = '  +2.89%  '

$ws.Range('E19').Value = '  +3.83%  '

$ws.Range('E20').Value = '  -1.51%  '

$ws.Range('D21').Value = "'75.51"
$ws.Range('E21').Value = '  +4.60%  '

$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = "'2.57"
$ws.Range('E22').Value = '  +11.22%  '

$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').Value = "'3.44"
$ws.Range('E23').Value = '  -0.41%  '

$ws.Range('D24').Value = "'256.94"
$ws.Range('E24').Value = '  +12.58%  '

$ws.Range('D25').Value = "'9.24"
$ws.Range('E25').Value = '  +1.72%  '

$ws.Range('E26').Value = '  +4.09%  '

$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('D28').Value = "'39.22"
$ws.Range('E28').Value = '  +2.82%  '

$ws.Range('E29').Value = '  +1.02%  '

$ws.Range('D30').Value = "'22.49"
$ws.Range('E30').Value = '  +7.39%  '

$ws.Range('D31').Value = "'173.64"
$ws.Range('E31').Value = '  -0.12%  '

$ws.Range('D32').Value = "'3.19"
$ws.Range('E32').Value = '  -0.24%  '

$ws.Range('E33').Value = '  +3.80%  '

$ws.Range('D34').Value = "'6.02"
$ws.Range('E34').Value = '  +8.00%  '

$ws.Range('E35').Value = '  +5.98%  '

$ws.Range('D36').Value = "'5.00"
$ws.Range('E36').Value = '  -1.22%  '

$ws.Range('D37').Value = "'4.16"
$ws.Range('E37').Value = '  -3.74%  '

$ws.Range('D38').Value = "'0.0376"
$ws.Range('E38').Value = '  -1.19%  '

$ws.Range('D39').Value = "'0.104"
$ws.Range('E39').Value = '  +1.87%  '

$ws.Range('D40').Value = "'2.73"
$ws.Range('E40').Value = '  +13.82%  '

$ws.Range('D41').Value = "'72.45"
$ws.Range('E41').Value = '  +0.92%  '

$ws.Range('D42').Value = "'1.49"
$ws.Range('E42').Value = '  +13.98%  '

$ws.Range('D43').Value = "'0.234"
$ws.Range('E43').Value = '  +1.52%  '

$ws.Range('E44').Value = '  +2.23%  '

$ws.Range('E45').Value = '  +0.18%  '

$ws.Range('E46').Value = '  +4.66%  '

$ws.Range('D47').Value = "'9.33"
$ws.Range('E47').Value = '  +11.09%  '

$ws.Range('D48').Value = "'110.99"
$ws.Range('E48').Value = '  +7.73%  '

$ws.Range('D49').Value = "'1.31"
$ws.Range('E49').Value = '  +0.59%  '

$ws.Range('E50').Value = '  +3.12%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = "'0.464"
$ws.Range('E51').Value = '  +6.60%  '
